# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list on Thu Mar 21 17:32:23 UTC 2024 with GitHub Actions".
#
# The sheet stores Price (D) and Volume(1h) (E) as plain text cells (they
# include thousands-separator dots and percent signs/padding), so several
# of the new Price values look numeric (e.g. "0.997", "12.00", "0.0000278")
# and would otherwise be silently re-interpreted/reformatted by Excel when
# assigned through .Value. To keep them as literal text we force those
# cells to the Text number format ("@") before writing the new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-converted to a number
# by Excel (losing the exact display text) -- force them to Text first.
$textCells = @(
    "D5",
    "D6",
    "D7",
    "D10",
    "D11",
    "D12",
    "D13",
    "D14",
    "D17",
    "D18",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D39",
    "D41",
    "D45",
    "D46",
    "D47",
    "D48",
    "D50",
    "D51"
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Row 2: Bitcoin
$ws.Range("D2").Value = '66.271.61'
$ws.Range("E2").Value = '  +3.26%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '3.492.39'
$ws.Range("E3").Value = '  +5.51%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.28%  '

# Row 5: BNB
$ws.Range("D5").Value = '559.35'
$ws.Range("E5").Value = '  +6.31%  '

# Row 6: Solana
$ws.Range("D6").Value = '183.27'
$ws.Range("E6").Value = '  +6.54%  '

# Row 7: XRP
$ws.Range("D7").Value = '0.637'
$ws.Range("E7").Value = '  +8.95%  '

# Row 8: LidoStakedEther
$ws.Range("D8").Value = '3.487.39'
$ws.Range("E8").Value = '  +5.49%  '

# Row 9: USDC
$ws.Range("E9").Value = '  +0.04%  '

# Row 10: Cardano
$ws.Range("D10").Value = '0.635'
$ws.Range("E10").Value = '  +5.44%  '

# Row 11: Dogecoin
$ws.Range("D11").Value = '0.154'
$ws.Range("E11").Value = '  +15.90%  '

# Row 12: Avalanche
$ws.Range("D12").Value = '54.97'
$ws.Range("E12").Value = '  +4.32%  '

# Row 13: ShibaInu
$ws.Range("D13").Value = '0.0000278'
$ws.Range("E13").Value = '  +8.16%  '

# Row 14: Polkadot
$ws.Range("D14").Value = '9.32'
$ws.Range("E14").Value = '  +4.58%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '4.081.09'
$ws.Range("E15").Value = '  +6.48%  '

# Row 16: WrappedEther
$ws.Range("D16").Value = '3.518.39'
$ws.Range("E16").Value = '  +6.45%  '

# Row 17: Chainlink
$ws.Range("D17").Value = '18.50'
$ws.Range("E17").Value = '  +6.39%  '

# Row 18: TRON
$ws.Range("D18").Value = '0.121'
$ws.Range("E18").Value = '  +3.98%  '

# Row 19: WrappedBTC
$ws.Range("D19").Value = '66.408.87'
$ws.Range("E19").Value = '  +3.69%  '

# Row 20: Uniswap
$ws.Range("D20").Value = '12.00'
$ws.Range("E20").Value = '  +7.55%  '

# Row 21: Polygon
$ws.Range("D21").Value = '0.997'
$ws.Range("E21").Value = '  +4.67%  '

# Row 22: BitcoinCash
$ws.Range("D22").Value = '417.58'
$ws.Range("E22").Value = '  +10.47%  '

# Row 23: PancakeSwap
$ws.Range("D23").Value = '4.09'
$ws.Range("E23").Value = '  +10.88%  '

# Row 24: Litecoin
$ws.Range("D24").Value = '85.98'
$ws.Range("E24").Value = '  +6.05%  '

# Row 25: Toncoin
$ws.Range("D25").Value = '4.28'
$ws.Range("E25").Value = '  +2.80%  '

# Row 26: ImmutableX
$ws.Range("D26").Value = '2.92'
$ws.Range("E26").Value = '  +8.17%  '

# Row 27: RenderToken
$ws.Range("D27").Value = '10.92'
$ws.Range("E27").Value = '  -1.45%  '

# Row 28: InternetComputer(DFINITY)
$ws.Range("D28").Value = '12.32'
$ws.Range("E28").Value = '  +9.90%  '

# Row 29: LEO
$ws.Range("D29").Value = '6.09'
$ws.Range("E29").Value = '  -1.21%  '

# Row 30: Filecoin
$ws.Range("D30").Value = '9.06'
$ws.Range("E30").Value = '  +11.82%  '

# Row 31: EthereumClassic
$ws.Range("D31").Value = '30.19'
$ws.Range("E31").Value = '  +5.43%  '

# Row 32: NEARProtocol
$ws.Range("D32").Value = '6.80'
$ws.Range("E32").Value = '  +3.83%  '

# Row 33: Bittensor
$ws.Range("D33").Value = '626.83'
$ws.Range("E33").Value = '  +0.10%  '

# Row 34: Cosmos
$ws.Range("D34").Value = '11.79'
$ws.Range("E34").Value = '  +5.72%  '

# Row 35: Hedera
$ws.Range("D35").Value = '0.111'
$ws.Range("E35").Value = '  +6.25%  '

# Row 36: OKB
$ws.Range("D36").Value = '60.13'
$ws.Range("E36").Value = '  +5.68%  '

# Row 37: Kaspa
$ws.Range("D37").Value = '0.148'
$ws.Range("E37").Value = '  +19.02%  '

# Row 38: PEPE
$ws.Range("D38").Value = '0.0₃0808'
$ws.Range("E38").Value = '  +9.01%  '

# Row 39: InjectiveProtocol
$ws.Range("D39").Value = '37.88'
$ws.Range("E39").Value = '  +5.48%  '

# Row 40: Dai
$ws.Range("E40").Value = '  -0.10%  '

# Row 41: TheGraph
$ws.Range("D41").Value = '0.384'
$ws.Range("E41").Value = '  +1.91%  '

# Row 42: Stacks
$ws.Range("E42").Value = '  +5.75%  '

# Row 43: Maker
$ws.Range("D43").Value = '3.086.17'
$ws.Range("E43").Value = '  +7.38%  '

# Row 44: FirstDigitalUSD
$ws.Range("E44").Value = '  +0.02%  '

# Row 45: Fetch.AI
$ws.Range("D45").Value = '2.60'
$ws.Range("E45").Value = '  +0.50%  '

# Row 46: ThetaToken
$ws.Range("D46").Value = '2.85'
$ws.Range("E46").Value = '  +9.48%  '

# Row 47: ApeXProtocol (was VeChain)
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '3.29'
$ws.Range("E47").Value = '  +8.24%  '

# Row 48: VeChain (was ApeXProtocol)
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = '0.0417'
$ws.Range("E48").Value = '  +5.37%  '

# Row 49: WEMIXToken
$ws.Range("E49").Value = '  +1.73%  '

# Row 50: Stellar
$ws.Range("D50").Value = '0.132'
$ws.Range("E50").Value = '  +6.56%  '

# Row 51: Monero
$ws.Range("D51").Value = '139.46'
$ws.Range("E51").Value = '  +1.49%  '
